$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 105, shifting existing rows 105:204 down to 106:205
$ws.Rows.Item(105).Insert()

# Populate the newly inserted row 105 with a new daily record
$ws.Range("A105").Value = 10
$ws.Range("B105").Value = "Vega Modelo de Temuco"
$ws.Range("C105").Value = "La Araucanía"
$ws.Range("D105").Value = 44589
$ws.Range("E105").Value = 9
$ws.Range("F105").Value = 100112039
$ws.Range("G105").Value = "Ciboulette"
$ws.Range("H105").Value = "Sin especificar"
$ws.Range("I105").Value = "Primera"
$ws.Range("J105").Value = 125
$ws.Range("K105").Value = 5000
$ws.Range("L105").Value = 5000
$ws.Range("M105").Value = 5000
$ws.Range("N105").Value = "$/docena de atados"
$ws.Range("O105").Value = "Provincia de Cautín"
$ws.Range("P105").Value = 1667
$ws.Range("Q105").Value = 3
$ws.Range("R105").Value = "Hortaliza"
